# Regenerate merged AHB files
# - Rename header row columns: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
# - Freeze the header row (row 1)
# - Wrap the data range in an Excel Table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header labels (row 1) ---------------------------------------
$oldToFV2410 = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
}

$newToFV2504 = @{
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $oldToFV2410.Keys) {
    $ws.Range($addr).Value = $oldToFV2410[$addr]
}

foreach ($addr in $newToFV2504.Keys) {
    $ws.Range($addr).Value = $newToFV2504[$addr]
}

# K1 ("diff") is unchanged.

# --- Freeze the header row ------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn the used range into a proper Excel Table (ListObject) ---------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), $null, 1)
$tbl.Name = "Table1"
